# Rename the sheet "Explorer_Election_16052022_0709" to "Sheet1"
# and move the current selection to cell E6, matching the saved
# workbook/worksheet state captured in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

$ws.Range("E6").Select()
